$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the comma-decimal header labels ("0,1" .. "0,6") with
# period-decimal versions ("0.1" .. "0.6"), keeping them as TEXT
# (shared-string) cells with their original number format/style.
# Temporarily switching to a text format while writing the value stops
# Excel from auto-converting the text to a real number, then restoring
# the original number format reuses the existing style (no new style
# record is left attached to the cells).
$origFormats = @{}
foreach ($addr in @("C1", "D1", "E1", "F1", "G1", "H1")) {
    $origFormats[$addr] = $ws.Range($addr).NumberFormat
}

$newValues = @{
    "C1" = "0.1"
    "D1" = "0.2"
    "E1" = "0.3"
    "F1" = "0.4"
    "G1" = "0.5"
    "H1" = "0.6"
}

foreach ($addr in @("C1", "D1", "E1", "F1", "G1", "H1")) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $newValues[$addr]
    $ws.Range($addr).NumberFormat = $origFormats[$addr]
}

$ws.Range("L10").Select()
